{"js": "// This document is a \"three-digit number \u00d7 one-digit number\" worksheet:\n// a table of math-problem cells such as \"321\u00d72=\". The commit regenerates\n// the worksheet with a new set of random problems, so every problem\n// <w:t> run's text is replaced by a new one (old -> new is a unique,\n// 1:1 mapping with no collisions between any \"old\" and any \"new\" value).\n// The title/date paragraph (\"2025-11-13 Thursday\") is left untouched.\nconst replacements = [\n  [\"321\u00d72=\", \"475\u00d75=\"],\n  [\"511\u00d75=\", \"959\u00d79=\"],\n  [\"340\u00d75=\", \"904\u00d76=\"],\n  [\"938\u00d73=\", \"394\u00d79=\"],\n  [\"944\u00d78=\", \"690\u00d74=\"],\n  [\"782\u00d76=\", \"591\u00d79=\"],\n  [\"899\u00d78=\", \"311\u00d76=\"],\n  [\"966\u00d77=\", \"152\u00d77=\"],\n  [\"907\u00d75=\", \"942\u00d79=\"],\n  [\"655\u00d75=\", \"756\u00d74=\"],\n  [\"537\u00d77=\", \"571\u00d77=\"],\n  [\"466\u00d74=\", \"220\u00d72=\"],\n  [\"426\u00d75=\", \"544\u00d72=\"],\n  [\"401\u00d79=\", \"283\u00d78=\"],\n  [\"201\u00d79=\", \"784\u00d73=\"],\n  [\"559\u00d77=\", \"297\u00d72=\"],\n  [\"940\u00d78=\", \"659\u00d72=\"],\n  [\"322\u00d75=\", \"999\u00d75=\"],\n  [\"623\u00d79=\", \"435\u00d76=\"],\n  [\"618\u00d77=\", \"205\u00d76=\"],\n  [\"296\u00d79=\", \"804\u00d76=\"],\n  [\"885\u00d77=\", \"528\u00d76=\"],\n  [\"114\u00d77=\", \"292\u00d74=\"],\n  [\"965\u00d78=\", \"128\u00d76=\"],\n  [\"277\u00d79=\", \"164\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# This document is a \"three-digit number x one-digit number\" worksheet:\n# a table of math-problem cells such as \"321x2=\". The commit regenerates\n# the worksheet with a new set of random problems, so every problem\n# cell's text is replaced by a new one (old -> new is a unique, 1:1\n# mapping with no collisions between any \"old\" and any \"new\" value).\n# The title/date paragraph (\"2025-11-13 Thursday\") is left untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"321\u00d72=\", \"475\u00d75=\"),\n    @(\"511\u00d75=\", \"959\u00d79=\"),\n    @(\"340\u00d75=\", \"904\u00d76=\"),\n    @(\"938\u00d73=\", \"394\u00d79=\"),\n    @(\"944\u00d78=\", \"690\u00d74=\"),\n    @(\"782\u00d76=\", \"591\u00d79=\"),\n    @(\"899\u00d78=\", \"311\u00d76=\"),\n    @(\"966\u00d77=\", \"152\u00d77=\"),\n    @(\"907\u00d75=\", \"942\u00d79=\"),\n    @(\"655\u00d75=\", \"756\u00d74=\"),\n    @(\"537\u00d77=\", \"571\u00d77=\"),\n    @(\"466\u00d74=\", \"220\u00d72=\"),\n    @(\"426\u00d75=\", \"544\u00d72=\"),\n    @(\"401\u00d79=\", \"283\u00d78=\"),\n    @(\"201\u00d79=\", \"784\u00d73=\"),\n    @(\"559\u00d77=\", \"297\u00d72=\"),\n    @(\"940\u00d78=\", \"659\u00d72=\"),\n    @(\"322\u00d75=\", \"999\u00d75=\"),\n    @(\"623\u00d79=\", \"435\u00d76=\"),\n    @(\"618\u00d77=\", \"205\u00d76=\"),\n    @(\"296\u00d79=\", \"804\u00d76=\"),\n    @(\"885\u00d77=\", \"528\u00d76=\"),\n    @(\"114\u00d77=\", \"292\u00d74=\"),\n    @(\"965\u00d78=\", \"128\u00d76=\"),\n    @(\"277\u00d79=\", \"164\u00d79=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 -> wdFindContinue, Replace:=2 -> wdReplaceAll\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
